$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title shape: merge "A" + " " + "slide" runs into a single run "A slide".
# A same-value (no-op) assignment would leave the multi-run structure
# untouched, and a simple substring edit would be treated as a minimal
# patch that also preserves the run split, so first push the text to an
# unrelated value, forcing a full rebuild into one run, then set the
# real target text.
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "placeholder"
$titleRange.Text = "A slide"

# Caption textbox: merge the many word/space runs into a single run
# "Just an image on this side" using the same two-step trick.
$captionRange = $s.Shapes.Item(4).TextFrame.TextRange
$captionRange.Text = "placeholder"
$captionRange.Text = "Just an image on this side"
